# Fix animation in ppt
#
# Slide 12 ("Docker-compose") has a bulleted content placeholder
# (shape id=5) with three click-triggered "Fade" entrance effects in
# its main animation sequence. The effects were wired to the wrong
# paragraphs (3rd bullet first, then 1st bullet, then 2nd bullet)
# instead of following reading order (1st, 2nd, then 3rd/4th
# paragraph). Re-point each effect at the correct paragraph so the
# click order matches the text order.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(12)
$seq = $s.TimeLine.MainSequence

# Effect.Paragraph is 1-based (paragraph 1 == a:p index 0 == p:pRg st="0").
$seq.Item(1).Paragraph = 1   # was targeting paragraph 4 (p:pRg st="3") -> now paragraph 1 (st="0")
$seq.Item(2).Paragraph = 2   # was targeting paragraph 1 (p:pRg st="0") -> now paragraph 2 (st="1")
$seq.Item(3).Paragraph = 4   # was targeting paragraph 2 (p:pRg st="1") -> now paragraph 4 (st="3")
